# Repull data, push all data, mean calculation
# Updates the dSF (column F) values for several rows with newly pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -1
$ws.Range("F6").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("F17").Value = -2
$ws.Range("F22").Value = 3
$ws.Range("F29").Value = 3
$ws.Range("F30").Value = -3
$ws.Range("F33").Value = 1
$ws.Range("F36").Value = 4
$ws.Range("F39").Value = -1
$ws.Range("F47").Value = 0
